$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("其他")
$ws.Rows(73).Insert()
